# Cotações atualizadas - 2025-10-13
# Adds the new daily quote row (row 39) to the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 39

# Date column (A) - keep same numeric/date style as the rows above it.
$ws.Cells.Item($newRow, 1).Value = 45943
$ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Quote columns (B:E) - stored as text, matching the existing rows.
$ws.Cells.Item($newRow, 2).Value = "21,5075"
$ws.Cells.Item($newRow, 3).Value = "15,4423"
$ws.Cells.Item($newRow, 4).Value = "15,4423"
$ws.Cells.Item($newRow, 5).Value = "15,4423"
